# Auto-generated script to apply scheduled-runner value updates
# to the Mandragora_Profits workbook (FFXIV leve profit tracker).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1525.4546
$ws.Range("J17").Value = 1525.4546
$ws.Range("L17").Value = 4576.3638
$ws.Range("N17").Value = -4912.3638
# Row 86
$ws.Range("H86").Value = 69236.28
$ws.Range("I86").Value = 112486.63
$ws.Range("J86").Value = 1271.4286
$ws.Range("K86").Value = 112486.63
$ws.Range("L86").Value = 1271.4286
$ws.Range("M86").Value = -111363.63
$ws.Range("N86").Value = -3517.4286
# Row 89
$ws.Range("H89").Value = 69236.28
$ws.Range("I89").Value = 112486.63
$ws.Range("J89").Value = 1271.4286
$ws.Range("K89").Value = 562433.15
$ws.Range("L89").Value = 6357.143
$ws.Range("M89").Value = -556817.15
$ws.Range("N89").Value = -17589.143
# Row 103
$ws.Range("H103").Value = 2974.5
$ws.Range("I103").Value = 199
$ws.Range("J103").Value = 5750
$ws.Range("K103").Value = 597
$ws.Range("L103").Value = 17250
$ws.Range("M103").Value = -11
$ws.Range("N103").Value = -18422
# Row 106
$ws.Range("H106").Value = 2666.3333
$ws.Range("I106").Value = 2467.6
$ws.Range("K106").Value = 2467.6
$ws.Range("M106").Value = -1836.6
# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
# Row 135
$ws.Range("H135").Value = 564.8444
$ws.Range("I135").Value = 390.92307
$ws.Range("J135").Value = 1695.3334
$ws.Range("K135").Value = 3518.30763
$ws.Range("L135").Value = 15258.0006
$ws.Range("M135").Value = -983.3076299999998
$ws.Range("N135").Value = -20328.0006

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7291.68
$ws.Range("I32").Value = 4895.435
$ws.Range("J32").Value = 20870.4
$ws.Range("K32").Value = 4895.435
$ws.Range("L32").Value = 20870.4
$ws.Range("M32").Value = -4608.435
$ws.Range("N32").Value = -21444.4
# Row 37
$ws.Range("H37").Value = 17206.8
$ws.Range("I37").Value = 17206.8
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 17206.8
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -16933.8
$ws.Range("N37").ClearContents()
# Row 61
$ws.Range("H61").Value = 1771.3096
$ws.Range("I61").Value = 1571.6957
$ws.Range("J61").Value = 2012.9474
$ws.Range("K61").Value = 1571.6957
$ws.Range("L61").Value = 2012.9474
$ws.Range("M61").Value = -1359.6957
$ws.Range("N61").Value = -2436.9474
# Row 74
$ws.Range("H74").Value = 1274.3188
$ws.Range("I74").Value = 895.7347
$ws.Range("J74").Value = 2201.85
$ws.Range("K74").Value = 895.7347
$ws.Range("L74").Value = 2201.85
$ws.Range("M74").Value = -21.73469999999998
$ws.Range("N74").Value = -3949.85
# Row 77
$ws.Range("H77").Value = 1274.3188
$ws.Range("I77").Value = 895.7347
$ws.Range("J77").Value = 2201.85
$ws.Range("K77").Value = 4478.6735
$ws.Range("L77").Value = 11009.25
$ws.Range("M77").Value = -110.6734999999999
$ws.Range("N77").Value = -19745.25
# Row 95
$ws.Range("H95").Value = 19652.111
$ws.Range("J95").Value = 19652.111
$ws.Range("L95").Value = 19652.111
$ws.Range("N95").Value = -25144.111
# Row 136
$ws.Range("H136").Value = 1771.3096
$ws.Range("I136").Value = 1571.6957
$ws.Range("J136").Value = 2012.9474
$ws.Range("K136").Value = 4715.0871
$ws.Range("L136").Value = 6038.8422
$ws.Range("M136").Value = -2165.0871
$ws.Range("N136").Value = -11138.8422

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 57
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
# Row 74
$ws.Range("H74").Value = 21167.5
$ws.Range("J74").Value = 21167.5
$ws.Range("L74").Value = 21167.5
$ws.Range("N74").Value = -23039.5
# Row 77
$ws.Range("H77").Value = 21167.5
$ws.Range("J77").Value = 21167.5
$ws.Range("L77").Value = 63502.5
$ws.Range("N77").Value = -72862.5
# Row 136
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 35
$ws.Range("H35").Value = 4611.4287
$ws.Range("I35").Value = 4611.4287
$ws.Range("K35").Value = 4611.4287
$ws.Range("M35").Value = -4317.4287
# Row 58
$ws.Range("H58").Value = 4230.375
$ws.Range("I58").Value = 2596.6216
$ws.Range("J58").Value = 7411.8945
$ws.Range("K58").Value = 2596.6216
$ws.Range("L58").Value = 7411.8945
$ws.Range("M58").Value = -2393.6216
$ws.Range("N58").Value = -7817.8945
# Row 99
$ws.Range("H99").Value = 3142.1155
$ws.Range("I99").Value = 3034.0625
$ws.Range("J99").Value = 3315
$ws.Range("K99").Value = 3034.0625
$ws.Range("L99").Value = 3315
$ws.Range("M99").Value = -1536.0625
$ws.Range("N99").Value = -6311
# Row 126
$ws.Range("H126").Value = 3142.1155
$ws.Range("I126").Value = 3034.0625
$ws.Range("J126").Value = 3315
$ws.Range("K126").Value = 9102.1875
$ws.Range("L126").Value = 9945
$ws.Range("M126").Value = -6632.1875
$ws.Range("N126").Value = -14885
# Row 134
$ws.Range("H134").Value = 1439.5245
$ws.Range("I134").Value = 736.35895
$ws.Range("J134").Value = 2686.0454
$ws.Range("K134").Value = 2209.07685
$ws.Range("L134").Value = 8058.1362
$ws.Range("M134").Value = 325.9231499999996
$ws.Range("N134").Value = -13128.1362
# Row 136
$ws.Range("H136").Value = 4230.375
$ws.Range("I136").Value = 2596.6216
$ws.Range("J136").Value = 7411.8945
$ws.Range("K136").Value = 7789.864799999999
$ws.Range("L136").Value = 22235.6835
$ws.Range("M136").Value = -5239.864799999999
$ws.Range("N136").Value = -27335.6835

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 70
$ws.Range("H70").Value = 2622
$ws.Range("I70").Value = 1960
$ws.Range("J70").Value = 2787.5
$ws.Range("K70").Value = 5880
$ws.Range("L70").Value = 8362.5
$ws.Range("M70").Value = -5565
$ws.Range("N70").Value = -8992.5
# Row 73
$ws.Range("H73").Value = 2622
$ws.Range("I73").Value = 1960
$ws.Range("J73").Value = 2787.5
$ws.Range("K73").Value = 5880
$ws.Range("L73").Value = 8362.5
$ws.Range("M73").Value = -4788
$ws.Range("N73").Value = -10546.5
# Row 75
$ws.Range("H75").Value = 1975
$ws.Range("I75").Value = 960
$ws.Range("J75").Value = 3666.6667
$ws.Range("K75").Value = 2880
$ws.Range("L75").Value = 11000.0001
$ws.Range("M75").Value = -1882
$ws.Range("N75").Value = -12996.0001
# Row 78
$ws.Range("H78").Value = 1975
$ws.Range("I78").Value = 960
$ws.Range("J78").Value = 3666.6667
$ws.Range("K78").Value = 8640
$ws.Range("L78").Value = 33000.0003
$ws.Range("M78").Value = -3648
$ws.Range("N78").Value = -42984.0003
# Row 87
$ws.Range("H87").Value = 6122.5
$ws.Range("I87").Value = 1980
$ws.Range("J87").Value = 6714.2856
$ws.Range("K87").Value = 5940
$ws.Range("L87").Value = 20142.8568
$ws.Range("M87").Value = -4692
$ws.Range("N87").Value = -22638.8568
# Row 90
$ws.Range("H90").Value = 6122.5
$ws.Range("I90").Value = 1980
$ws.Range("J90").Value = 6714.2856
$ws.Range("K90").Value = 17820
$ws.Range("L90").Value = 60428.5704
$ws.Range("M90").Value = -11580
$ws.Range("N90").Value = -72908.5704
# Row 129
$ws.Range("H129").Value = 2505.0435
$ws.Range("I129").Value = 632.2222
$ws.Range("J129").Value = 3709
$ws.Range("K129").Value = 1896.6666
$ws.Range("L129").Value = 11127
$ws.Range("M129").Value = 3103.3334
$ws.Range("N129").Value = -21127
# Row 131
$ws.Range("H131").Value = 404093.84
$ws.Range("I131").Value = 455.55554
$ws.Range("J131").Value = 569218.5600000001
$ws.Range("K131").Value = 1366.66662
$ws.Range("L131").Value = 1707655.68
$ws.Range("M131").Value = 3673.33338
$ws.Range("N131").Value = -1717735.68

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 928414.7
$ws.Range("I132").Value = 1813327
$ws.Range("J132").Value = 3279
$ws.Range("K132").Value = 5439981
$ws.Range("L132").Value = 9837
$ws.Range("M132").Value = -5437451
$ws.Range("N132").Value = -14897

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 13335497
$ws.Range("I132").Value = 21741466
$ws.Range("J132").Value = 1890.8276
$ws.Range("K132").Value = 65224398
$ws.Range("L132").Value = 5672.4828
$ws.Range("M132").Value = -65221868
$ws.Range("N132").Value = -10732.4828
# Row 136
$ws.Range("H136").Value = 8929913
$ws.Range("I136").Value = 12821229
$ws.Range("J136").Value = 2776.4707
$ws.Range("K136").Value = 38463687
$ws.Range("L136").Value = 8329.4121
$ws.Range("M136").Value = -38461137
$ws.Range("N136").Value = -13429.4121

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1453.494
$ws.Range("I132").Value = 1333.0667
$ws.Range("J132").Value = 1596.1052
$ws.Range("K132").Value = 3999.2001
$ws.Range("L132").Value = 4788.3156
$ws.Range("M132").Value = -1469.2001
$ws.Range("N132").Value = -9848.3156
